$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "USA"
$ws.Range("C9").Value = "GBR"
$ws.Range("D9").Value = "GER"

$ws.Range("C15").Select()
